# Flow Order.xlsx update
#   - add 10 new "LULC" flow types (DOM/Live CO2/CO/CH4, Transfer DOM,
#     Harvest DOM/Live, Mortality Live) to the 'Flow Order' sheet
#   - re-sort the Flow Type/Order table (columns C:D) alphabetically by
#     Flow Type, which is what re-orders/interleaves the new rows among
#     the existing ones
#   - grow the AutoFilter / filter-database range to match the new data
#     extent (C1:D43 instead of A1:D33)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow Order")
$ws.Activate()

# Drop the existing AutoFilter so we can re-create it over the new range
# once the new rows are in place.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# New Flow Type / Order rows, appended below the current data (row 34
# onward). They get the plain/default cell style -- same as the rest of
# row 26-40 end up with after the sort below, since the style actually
# used for display on columns C/D comes from the column defaults
# (style ids 7 and 9 declared on <col>).
$newRows = @(
    @("LULC: Emission DOM CO2",  10),
    @("LULC: Emission DOM CO",   10),
    @("LULC: Emission DOM CH4",  10),
    @("LULC: Transfer DOM",      10),
    @("LULC: Harvest DOM",       10),
    @("LULC: Emission Live CO2", 11),
    @("LULC: Emission Live CO",  11),
    @("LULC: Emission Live CH4", 11),
    @("LULC: Harvest Live",      11),
    @("LULC: Mortality Live",    11)
)

$r = 34
foreach ($row in $newRows) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Style = "Normal"
    $cCell.Value2 = $row[0]

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Style = "Normal"
    $dCell.Value2 = $row[1]

    $r = $r + 1
}

# Sort the whole Flow Type/Order block (now rows 2-43) by the Flow Type
# column (C), ascending, header row included -- this interleaves the new
# LULC rows alphabetically among the pre-existing ones.
$sortRange = $ws.Range("C1:D43")
$sortKey = $ws.Range("C1:C43")
$sortRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)

# Re-apply the AutoFilter over the grown data range (C1:D43).
$ws.Range("C1:D43").AutoFilter() | Out-Null

# The hidden _FilterDatabase name needs to point at the same, grown range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Flow Order'!`$C`$1:`$D`$43"
    }
}

# Column C widened to fit the new (longer) Flow Type labels (best-fit,
# approximates the ~49.57-character width Excel computes for the longest
# new label, "LULC: Emission Live CO2").
$ws.Columns("C:C").ColumnWidth = 48.6

# Leave the selection where the edit ended.
$ws.Range("I36").Select() | Out-Null
